# Insert a new data row at row 624 (shifting the existing rows 624-673 down
# to 625-674), then populate the new row with the reported price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("624:624").Insert()

$ws.Range("A624").Value = 10
$ws.Range("B624").Value = "Vega Modelo de Temuco"
$ws.Range("C624").Value = "La Araucanía"
$ws.Range("D624").Value = 45265
$ws.Range("E624").Value = 9
$ws.Range("F624").Value = 100112024
$ws.Range("G624").Value = "Choclo"
$ws.Range("H624").Value = "Dulce o Americano"
$ws.Range("I624").Value = "Primera"
$ws.Range("J624").Value = 3000
$ws.Range("K624").Value = 700
$ws.Range("L624").Value = 700
$ws.Range("M624").Value = 700
$ws.Range("N624").Value = "`$/unidad"
$ws.Range("O624").Value = "Región de O'Higgins"
$ws.Range("P624").Value = 700
$ws.Range("Q624").Value = 1
$ws.Range("R624").Value = "Hortaliza"
